# Updated cryptos list with GitHub Actions
# Applies the latest price / 1h-volume-change snapshot to the cryptos sheet,
# and fixes the Monero / EthereumClassic row ordering (rows 31-32 swapped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bitcoin
$ws.Range("D2").Value = "42.998.47"
$ws.Range("E2").Value = "  +0.55%  "

# Ethereum
$ws.Range("D3").Value = "2.281.27"
$ws.Range("E3").Value = "  +1.58%  "

# TetherUSD
$ws.Range("E4").Value = "  +0.29%  "

# Solana
$ws.Range("D5").Value = "'112.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.47%  "

# BNB
$ws.Range("D6").Value = "'308.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.24%  "

# XRP
$ws.Range("E7").Value = "  -0.53%  "

# USDC
$ws.Range("E8").Value = "  +0.10%  "

# Cardano
$ws.Range("D9").Value = "'0.613"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.85%  "

# Avalanche
$ws.Range("D10").Value = "'44.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.15%  "

# Dogecoin
$ws.Range("D11").Value = "'0.0926"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.20%  "

# OKB
$ws.Range("E12").Value = "  -1.54%  "

# Polkadot
$ws.Range("D13").Value = "'8.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.74%  "

# Polygon
$ws.Range("D14").Value = "'1.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +21.48%  "

# TRON
$ws.Range("E15").Value = "  -0.09%  "

# Chainlink
$ws.Range("D16").Value = "'15.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.28%  "

# WrappedliquidstakedEther2.0
$ws.Range("D17").Value = "2.621.54"
$ws.Range("E17").Value = "  +1.39%  "

# WrappedEther
$ws.Range("D18").Value = "2.277.79"
$ws.Range("E18").Value = "  +1.48%  "

# WrappedBTC
$ws.Range("D19").Value = "42.930.51"
$ws.Range("E19").Value = "  +0.36%  "

# ShibaInu
$ws.Range("E20").Value = "  -0.66%  "

# Uniswap
$ws.Range("E21").Value = "  -3.58%  "

# Litecoin
$ws.Range("D22").Value = "'76.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.88%  "

# PancakeSwap
$ws.Range("D23").Value = "'3.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.42%  "

# ImmutableX
$ws.Range("E24").Value = "  +4.34%  "

# BitcoinCash
$ws.Range("D25").Value = "'255.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.49%  "

# InternetComputer(DFINITY)
$ws.Range("D26").Value = "'8.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.93%  "

# Cosmos
$ws.Range("D27").Value = "'11.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.99%  "

# Dai
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.06%  "

# Toncoin
$ws.Range("E29").Value = "  -0.43%  "

# InjectiveProtocol
$ws.Range("D30").Value = "'38.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.61%  "

# row 31 (Monero <-> EthereumClassic swap)
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'22.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.18%  "

# row 32 (EthereumClassic <-> Monero swap)
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'174.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.54%  "

# Hedera
$ws.Range("D34").Value = "'0.0900"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.65%  "

# Filecoin
$ws.Range("E35").Value = "  +0.87%  "

# RenderToken
$ws.Range("D36").Value = "'5.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.08%  "

# Stellar
$ws.Range("E37").Value = "  +0.11%  "

# NEARProtocol
$ws.Range("D38").Value = "'4.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.04%  "

# VeChain
$ws.Range("D39").Value = "'0.0377"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.52%  "

# Kaspa
$ws.Range("D40").Value = "'0.104"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.44%  "

# LidoDAOToken
$ws.Range("E41").Value = "  -3.68%  "

# MultiversX
$ws.Range("D42").Value = "'72.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.14%  "

# Algorand
$ws.Range("D43").Value = "'0.231"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.61%  "

# FirstDigitalUSD
$ws.Range("E44").Value = "  +0.22%  "

# Celestia
$ws.Range("D45").Value = "'12.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.76%  "

# ARBITRUM
$ws.Range("E46").Value = "  +1.86%  "

# THORChain
$ws.Range("D47").Value = "'5.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.33%  "

# Aave
$ws.Range("D48").Value = "'108.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.10%  "

# FraxShare
$ws.Range("E49").Value = "  +2.81%  "

# TrustWalletToken
$ws.Range("D50").Value = "'1.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.48%  "

# ordi
$ws.Range("D51").Value = "'71.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.00%  "
